$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.250.33"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.585.10"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "209.07"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "`'0.500"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "0.0611"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D13").Value = "1.616.50"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "26.249.48"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("D20").Value = "212.54"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "4.26"
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Value = "8.86"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "`'144.60"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "7.03"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "`'15.30"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "`'3.00"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "1.287.65"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").Value = "1.13"
$ws.Range("E38").Value = "  -7.36%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "62.59"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "0.762"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "1.719.63"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "88.69"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").Value = "0.0999"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "0.0₇0970"
$ws.Range("E51").Value = "  -8.34%  "
